$wb = $excel.ActiveWorkbook
$testData = $wb.Worksheets.Item("TestData")

# Fill in the missing test data value for the appium test case (cell D4).
# D4 already carried a specific number-format style before this edit; a
# plain Value assignment in this engine resets the cell's style index, so
# we restore the original formatting by pasting the format (only) from a
# neighboring cell that already uses that same style.
$testData.Range("D4").Value = "appium"
$testData.Range("C6").Copy()
$testData.Range("D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Switch to (activate) the "TestData" sheet and select the cell that was
# just edited - this becomes the new active tab/selection for the workbook,
# replacing "RunManager" as the active sheet.
$testData.Activate()
$testData.Range("D4").Select() | Out-Null
